# Table_2.xlsx — refresh recalculated incidence/prevalence figures in column C
# and size columns B/C to fit their (now wider) contents, per the
# "Fixed the incidence and prevalance code ... re-written the results table"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated numeric results (column C, rows 2-11)
$ws.Range("C2").Value  = 754892.98513708881
$ws.Range("C3").Value  = 142559.35349296677
$ws.Range("C4").Value  = 71248.97183725526
$ws.Range("C5").Value  = 541084.6598068655
$ws.Range("C6").Value  = 28039.031803541533
$ws.Range("C7").Value  = 134207.56925177248
$ws.Range("C8").Value  = 190481.30786526375
$ws.Range("C9").Value  = 223521.99651547629
$ws.Range("C10").Value = 178497.47097301373
$ws.Range("C11").Value = 145.60872802009752

# Column widths now explicitly sized to fit the content (bestFit in Excel)
$ws.Range("B:B").ColumnWidth = 14.83
$ws.Range("C:C").ColumnWidth = 11.1
